# Update countries & provincias Spain
#
# The upstream COVID-19 country feed refreshed between the 13:50 and 14:20
# snapshots: several countries' case counts changed, and a handful of rows
# swap places because the sheet is kept sorted descending by 'Casos totales'
# (column B). This script rewrites the timestamp banner and every row whose
# data actually changed, landing on the exact post-refresh values/order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Banner timestamp (cell A1)
$ws.Range("A1").Value = "Datos actualizados a 2 de Abril de 2020 a las 14:20"

# Row data: row number, Pais, Casos totales, Nuevos casos, Casos activos,
# Recuperados, Casos criticos, Muertes hoy, Muertes
$rows = @(
    @(15, "Paises Bajos", 14697, 1083, 250, 13108, 1053, 166, 1339),
    @(19, "Portugal", 9034, 783, 68, 8757, 230, 22, 209),
    @(20, "Brasil", 6932, 52, 127, 6559, 296, 4, 246),
    @(22, "Suecia", 5466, 519, 103, 5081, 429, 43, 282),
    @(23, "Australia", 5137, 89, 345, 4767, 50, 2, 25),
    @(24, "Noruega", 5071, 194, 13, 5012, 105, 2, 46),
    @(38, "India", 2032, 34, 150, 1824, 0, 0, 58),
    @(53, "Croacia", 1011, 48, 88, 916, 34, 1, 7),
    @(54, "Singapur", 1000, 0, 245, 751, 24, 1, 4),
    @(66, "Marruecos", 676, 22, 29, 607, 1, 1, 40),
    @(67, "Armenia", 663, 92, 33, 625, 30, 1, 5),
    @(68, "Lituania", 649, 68, 7, 633, 11, 1, 9),
    @(71, "Bosnia y Herzegovina", 518, 59, 20, 483, 4, 2, 15),
    @(79, "Kazajistan", 423, 43, 27, 393, 6, 0, 3),
    @(81, "Republica de Macedonia", 384, 30, 17, 356, 4, 0, 11),
    @(82, "Costa Rica", 375, 0, 4, 369, 9, 0, 2),
    @(87, "Camerun", 284, 51, 10, 267, 0, 1, 7),
    @(88, "Burkina Faso", 282, 0, 46, 220, 0, 0, 16),
    @(89, "Reunion", 281, 0, 40, 241, 3, 0, 0),
    @(90, "Jordania", 278, 0, 36, 237, 5, 0, 5),
    @(91, "Albania", 277, 18, 67, 194, 7, 1, 16)
)

foreach ($r in $rows) {
    $rowNum = $r[0]
    for ($col = 1; $col -lt $r.Length; $col++) {
        $ws.Cells.Item($rowNum, $col).Value = $r[$col]
    }
}
